{"js": "// Rewrite the \"Add the following line to it...\" block in the plotting\n// instructions section so that:\n//   1. The numbered list item becomes \"Enter a line with the following format: \"\n//   2. The italicized format line \"<save_name>  <repetitons>  <BehaviorSpace argument>\"\n//      moves up to sit directly beneath it.\n//   3. The old \"Tip: Copy/Paste...\" line becomes a \"Copy/Paste\" (bold) led\n//      sentence with extra guidance text.\n//   4. The \"Note: ...\" line is reworded slightly and the \"_GoBack\" bookmark\n//      (previously inside the numbered list item) is relocated into it.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the four paragraphs we need by their distinctive original text so\n// the script is resilient to any surrounding content shifting around.\nlet idxAdd = -1, idxTip = -1, idxNote = -1, idxFormat = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (idxAdd === -1 && t.indexOf(\"Add the following line to it\") !== -1) {\n    idxAdd = i;\n  } else if (idxTip === -1 && t.indexOf(\"Tip: Copy/Paste the last argument\") !== -1) {\n    idxTip = i;\n  } else if (idxNote === -1 && t.indexOf(\"Note: The format for plotting parameters\") !== -1) {\n    idxNote = i;\n  } else if (idxFormat === -1 && t.indexOf(\"<save_name>\") !== -1) {\n    idxFormat = i;\n  }\n}\n\nif (idxAdd === -1 || idxTip === -1 || idxNote === -1 || idxFormat === -1) {\n  throw new Error(\n    \"Could not locate all four target paragraphs: \" +\n    JSON.stringify({ idxAdd, idxTip, idxNote, idxFormat })\n  );\n}\n\nconst pAdd = paragraphs.items[idxAdd];\nconst pTip = paragraphs.items[idxTip];\nconst pNote = paragraphs.items[idxNote];\nconst pFormat = paragraphs.items[idxFormat];\n\n// 1) Rewrite the numbered list paragraph's text (this also removes the\n//    \"_GoBack\" bookmark that used to live inside it).\npAdd.clear();\nawait context.sync();\npAdd.insertText(\"Enter a line with the following format\", Word.InsertLocation.end);\nawait context.sync();\npAdd.insertText(\": \", Word.InsertLocation.end);\nawait context.sync();\n\n// 2) Insert a fresh copy of the italic \"<save_name> ...\" line directly\n//    before the \"Tip\" paragraph (inserting \"before\" a non-list paragraph\n//    keeps the same ListParagraph/indent/spacing formatting without\n//    picking up numbering).\nconst newFormatPara = pTip.insertParagraph(\n  \"<save_name>  <repetitons>  <BehaviorSpace argument>\",\n  Word.InsertLocation.before\n);\nnewFormatPara.font.italic = true;\nawait context.sync();\n\n// Remove the now-duplicated original italic paragraph.\npFormat.delete();\nawait context.sync();\n\n// 3) Rework the \"Tip\" paragraph into the new \"Copy/Paste ...\" sentence.\npTip.clear();\nawait context.sync();\nconst copyPasteRun = pTip.insertText(\"Copy/Paste\", Word.InsertLocation.end);\ncopyPasteRun.font.bold = true;\nawait context.sync();\nconst restRun = pTip.insertText(\n  \" the last argument from BehaviorSpace. Typing it may produce errors\",\n  Word.InsertLocation.end\n);\nrestRun.font.bold = false;\nawait context.sync();\n\n// 4) Reword the \"Note\" paragraph and move the \"_GoBack\" bookmark into it,\n//    sitting right after \"Note: The\".\npNote.clear();\nawait context.sync();\npNote.insertText(\n  \"Note: The format for plotting parameters will always be the same.\",\n  Word.InsertLocation.end\n);\nawait context.sync();\n\nconst searchResults = pNote.search(\"The\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\nconst theRange = searchResults.items[0].getRange(Word.RangeLocation.after);\nawait context.sync();\ntheRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Rewrite the \"Add the following line to it...\" block in the plotting\n# instructions section so that:\n#   1. The numbered list item becomes \"Enter a line with the following format: \"\n#   2. The italicized format line \"<save_name>  <repetitons>  <BehaviorSpace argument>\"\n#      moves up to sit directly beneath it.\n#   3. The old \"Tip: Copy/Paste...\" line becomes a \"Copy/Paste\" (bold) led\n#      sentence with extra guidance text.\n#   4. The \"Note: ...\" line is reworded slightly and the \"_GoBack\" bookmark\n#      (previously inside the numbered list item) is relocated into it.\n\n$d = $word.ActiveDocument\n\n# --- locate the four target paragraphs by their distinctive text -----------\n$idxAdd = -1\n$idxTip = -1\n$idxNote = -1\n$idxFormat = -1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($idxAdd -eq -1 -and $t -like \"*Add the following line to it*\") {\n        $idxAdd = $p.Index\n    } elseif ($idxTip -eq -1 -and $t -like \"*Tip: Copy/Paste the last argument*\") {\n        $idxTip = $p.Index\n    } elseif ($idxNote -eq -1 -and $t -like \"*Note: The format for plotting parameters*\") {\n        $idxNote = $p.Index\n    } elseif ($idxFormat -eq -1 -and $t -like \"*<save_name>*\") {\n        $idxFormat = $p.Index\n    }\n}\n\n# 1) Rewrite the numbered list paragraph's text. Using Find/Replace (rather\n#    than Range.Text=) so every run inside the paragraph (and the old\n#    \"_GoBack\" bookmark that lived inside it) is cleanly replaced.\n$addRange = $d.Paragraphs($idxAdd).Range\n$addFind = $addRange.Find\n$addFind.Execute( `\n    \"Add the following line to it*1]\", `\n    $false, $false, $true, $false, $false, $true, 1, $false, `\n    \"Enter a line with the following format: \", 2)\n\n# 2) Insert a fresh copy of the italic \"<save_name> ...\" line directly\n#    before the \"Tip\" paragraph (inserting \"before\" a non-list paragraph\n#    keeps the same ListParagraph/indent/spacing formatting without\n#    picking up numbering).\n$tipRange = $d.Paragraphs($idxTip).Range\n$tipRange.InsertParagraphBefore()\n\n# The blank paragraph just inserted is now at $idxTip; \"Tip\" shifted to\n# $idxTip + 1, \"Note\" to $idxNote + 1, and the original format line to\n# $idxFormat + 1.\n$newFormatPara = $d.Paragraphs($idxTip)\n$newFormatRange = $newFormatPara.Range\n$newFormatRange.InsertAfter(\"<save_name>  <repetitons>  <BehaviorSpace argument>\")\n$newFormatRange2 = $d.Paragraphs($idxTip).Range\n$newFormatRange2.Font.Italic = 1\n\n$idxTip = $idxTip + 1\n$idxNote = $idxNote + 1\n$idxFormat = $idxFormat + 1\n\n# Remove the now-duplicated original italic paragraph.\n$d.Paragraphs($idxFormat).Range.Delete()\n\n# 3) Rework the \"Tip\" paragraph into the new \"Copy/Paste ...\" sentence,\n#    applying Bold only to the \"Copy/Paste\" word.\n$tipClearRange = $d.Paragraphs($idxTip).Range\n$tipClearFind = $tipClearRange.Find\n$tipClearFind.Execute( `\n    \"Tip: Copy/Paste the last argument from BehaviorSpace\", `\n    $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n\n$tipStart = $d.Paragraphs($idxTip).Range.Start\n$d.Paragraphs($idxTip).Range.InsertAfter(\"Copy/Paste\")\n$boldRange = $d.Range($tipStart, $tipStart + 10)\n$boldRange.Font.Bold = 1\n\n$restText = \" the last argument from BehaviorSpace. Typing it may produce errors\"\n$d.Paragraphs($idxTip).Range.InsertAfter($restText)\n$restStart = $tipStart + 10\n$restRange = $d.Range($restStart, $restStart + $restText.Length)\n$restRange.Font.Bold = 0\n\n# 4) Reword the \"Note\" paragraph and move the \"_GoBack\" bookmark into it,\n#    sitting right after \"Note: The\".\n$noteRange = $d.Paragraphs($idxNote).Range\n$noteFind = $noteRange.Find\n$noteFind.Execute( `\n    \"Note: The format for plotting parameters will always be the same.\", `\n    $false, $false, $false, $false, $false, $true, 1, $false, `\n    \"Note: The format for plotting parameters will always be the same.\", 2)\n\n$bookmarkRange = $d.Paragraphs($idxNote).Range\n$bookmarkFind = $bookmarkRange.Find\n$bookmarkFind.Execute(\"Note: The\")\n$bookmarkRange.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
